# "Our primary goal is to create a single landing page for an imaginary
# client. This site must behave ..." becomes "Our primary goal is to
# create a single landing page for a client. The site must behave ...".
#
# The surviving text ends up split across three runs (..."Th" / "e" /
# " site"...) exactly as the canonical OOXML diff shows, while the
# following run ("sheets. All tasks...") stays a separate run too.

$d = $word.ActiveDocument

# --- Step 1: drop "an imaginary " before "client" -----------------------
$content = $d.Content
$content.Find.Execute("an imaginary client", $true, $false, $false, $false, `
    $false, $true, 1, $false, "a client", 2) | Out-Null

# --- Step 2: "This site" -> "The site" (only the "is" -> "e" part) ------
$content = $d.Content
$bodyText = $content.Text
$idx = $bodyText.IndexOf("This site")
$start = $content.Start + $idx + 2   # right after "Th"
$end = $start + 2                     # the "is" that needs to become "e"
$mid = $d.Range($start, $end)
$mid.Text = "e"

# --- Step 3: re-create the run boundaries seen in the target XML --------
# Toggling (and un-toggling) a character property on a narrow range is
# enough to make the engine keep that range as its own <w:r> once the
# surrounding text has already been rewritten above, without altering
# the run's effective formatting.
function SplitRunBoundary($absStart, $absEnd) {
    $rr = $d.Range($absStart, $absEnd)
    $rr.Bold = 1
    $rr2 = $d.Range($absStart, $absEnd)
    $rr2.Bold = 0
}

# 3a: isolate the single "e" that replaced "is" into its own run
$content = $d.Content
$bodyText = $content.Text
$idxE = $bodyText.IndexOf("The site")
$eStart = $content.Start + $idxE + 2
$eEnd = $eStart + 1
SplitRunBoundary $eStart $eEnd

# 3b: restore the pre-existing run break right before "sheets." (the
# rewrite above re-flattens every identically-formatted run in the
# paragraph, so this boundary has to be put back too)
$content = $d.Content
$bodyText = $content.Text
$idxSheets = $bodyText.IndexOf("sheets.")
$tailMarker = "in a timely manner."
$idxTailEnd = $bodyText.IndexOf($tailMarker) + $tailMarker.Length
$sheetsStart = $content.Start + $idxSheets
$sheetsEnd = $content.Start + $idxTailEnd
SplitRunBoundary $sheetsStart $sheetsEnd
